$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 16-29 hold a table of "Valor Mora" (overdue amounts) per worker
# per period. The previous layout interleaved the two workers row by row;
# the new layout groups all rows for KARETH first (periods 2311 -> 2305,
# descending), followed by all rows for KELLY (periods 2311 -> 2305,
# descending). The old 37333 amount (tied to period 2311) now follows its
# period for both workers, while every other period uses 46400.

$rows = @(
    @{ Row = 16; Doc = "22802542";   Name = "KARETH LAVINIA GUZMAN PAJARO"; Period = "2311"; Valor = 37333 },
    @{ Row = 17; Doc = "22802542";   Name = "KARETH LAVINIA GUZMAN PAJARO"; Period = "2310"; Valor = 46400 },
    @{ Row = 18; Doc = "22802542";   Name = "KARETH LAVINIA GUZMAN PAJARO"; Period = "2309"; Valor = 46400 },
    @{ Row = 19; Doc = "22802542";   Name = "KARETH LAVINIA GUZMAN PAJARO"; Period = "2308"; Valor = 46400 },
    @{ Row = 20; Doc = "22802542";   Name = "KARETH LAVINIA GUZMAN PAJARO"; Period = "2307"; Valor = 46400 },
    @{ Row = 21; Doc = "22802542";   Name = "KARETH LAVINIA GUZMAN PAJARO"; Period = "2306"; Valor = 46400 },
    @{ Row = 22; Doc = "22802542";   Name = "KARETH LAVINIA GUZMAN PAJARO"; Period = "2305"; Valor = 46400 },
    @{ Row = 23; Doc = "1091680446"; Name = "KELLY DAYANA ASCANIO TORRES";  Period = "2311"; Valor = 37333 },
    @{ Row = 24; Doc = "1091680446"; Name = "KELLY DAYANA ASCANIO TORRES";  Period = "2310"; Valor = 46400 },
    @{ Row = 25; Doc = "1091680446"; Name = "KELLY DAYANA ASCANIO TORRES";  Period = "2309"; Valor = 46400 },
    @{ Row = 26; Doc = "1091680446"; Name = "KELLY DAYANA ASCANIO TORRES";  Period = "2308"; Valor = 46400 },
    @{ Row = 27; Doc = "1091680446"; Name = "KELLY DAYANA ASCANIO TORRES";  Period = "2307"; Valor = 46400 },
    @{ Row = 28; Doc = "1091680446"; Name = "KELLY DAYANA ASCANIO TORRES";  Period = "2306"; Valor = 46400 },
    @{ Row = 29; Doc = "1091680446"; Name = "KELLY DAYANA ASCANIO TORRES";  Period = "2305"; Valor = 46400 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc
    $ws.Cells.Item($r.Row, 4).Value = $r.Name
    $ws.Cells.Item($r.Row, 5).Value = $r.Period
    $ws.Cells.Item($r.Row, 6).Value = $r.Valor
}
